$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report gained a new "2022" column (J). Copy the formatting already
# used for the 2021 column (I) into the new J column for the data block
# (rows 3-14), then fill in the 2022 figures.
$ws.Range("I3:I14").Copy() | Out-Null
$ws.Range("J3:J14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("J4").Value = 2022
$ws.Range("J5").Value = 96.4
$ws.Range("J6").Value = 96.4
$ws.Range("J7").Value = 97.9
$ws.Range("J8").Value = 95.3
$ws.Range("J9").Value = 93.8
$ws.Range("J10").Value = 95.5
$ws.Range("J11").Value = 94.4
$ws.Range("J12").Value = 95
$ws.Range("J13").Value = 98.7
$ws.Range("J14").Value = 97.3

# Match the saved selection shown in the sheet view.
$ws.Range("L10").Select() | Out-Null
